$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 (Generation 0,1,2) -> 7586
$ws.Range("C2:C4").Value = 7586

# Rows 5-252 (Generation 3..250) -> 7310
$ws.Range("C5:C252").Value = 7310
